$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.076.11'
$ws.Range("E2").Value = '  -0.70%  '

# Row 3
$ws.Range("D3").Value = '2.015.90'
$ws.Range("E3").Value = '  -1.71%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = "'226.44"
$ws.Range("E5").Value = '  -1.72%  '

# Row 6
$ws.Range("E6").Value = '  -2.32%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").Value = "'54.76"
$ws.Range("E8").Value = '  -4.05%  '

# Row 9
$ws.Range("D9").Value = "'0.379"
$ws.Range("E9").Value = '  -1.37%  '

# Row 10
$ws.Range("E10").Value = '  +2.39%  '

# Row 11
$ws.Range("E11").Value = '  -3.14%  '

# Row 12
$ws.Range("D12").Value = '2.315.67'
$ws.Range("E12").Value = '  -1.64%  '

# Row 13
$ws.Range("E13").Value = '  -3.03%  '

# Row 14
$ws.Range("D14").Value = "'20.37"
$ws.Range("E14").Value = '  -0.88%  '

# Row 15
$ws.Range("D15").Value = "'0.741"
$ws.Range("E15").Value = '  -2.03%  '

# Row 16
$ws.Range("D16").Value = "'5.13"
$ws.Range("E16").Value = '  -2.13%  '

# Row 17
$ws.Range("D17").Value = '2.004.70'
$ws.Range("E17").Value = '  -2.09%  '

# Row 18
$ws.Range("D18").Value = '37.012.35'
$ws.Range("E18").Value = '  -0.77%  '

# Row 19
$ws.Range("D19").Value = "'6.14"
$ws.Range("E19").Value = '  +2.53%  '

# Row 20
$ws.Range("D20").Value = "'68.84"
$ws.Range("E20").Value = '  -1.38%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  -0.45%  '

# Row 22
$ws.Range("D22").Value = "'224.08"
$ws.Range("E22").Value = '  -1.15%  '

# Row 23
$ws.Range("E23").Value = '  -0.07%  '

# Row 24
$ws.Range("E24").Value = '  +1.46%  '

# Row 25
$ws.Range("D25").Value = "'2.20"
$ws.Range("E25").Value = '  -5.37%  '

# Row 26
$ws.Range("D26").Value = "'165.25"
$ws.Range("E26").Value = '  -2.11%  '

# Row 27
$ws.Range("D27").Value = "'9.17"
$ws.Range("E27").Value = '  -3.24%  '

# Row 28
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = "'1.36"
$ws.Range("E28").Value = '  +0.85%  '

# Row 29
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = "'0.126"
$ws.Range("E29").Value = '  -3.84%  '

# Row 30
$ws.Range("D30").Value = "'18.71"
$ws.Range("E30").Value = '  -2.28%  '

# Row 31
$ws.Range("E31").Value = '  -3.61%  '

# Row 32
$ws.Range("D32").Value = "'4.54"
$ws.Range("E32").Value = '  +0.41%  '

# Row 33
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = '  -1.36%  '

# Row 34
$ws.Range("D34").Value = "'4.41"
$ws.Range("E34").Value = '  -3.17%  '

# Row 35
$ws.Range("D35").Value = "'2.34"
$ws.Range("E35").Value = '  -5.79%  '

# Row 36
$ws.Range("E36").Value = '  +1.88%  '

# Row 37
$ws.Range("E37").Value = '  +0.11%  '

# Row 38
$ws.Range("D38").Value = "'3.14"
$ws.Range("E38").Value = '  -4.24%  '

# Row 39
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = '  +2.69%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "'0.0218"
$ws.Range("E40").Value = '  -3.83%  '

# Row 41
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.483.26'
$ws.Range("E41").Value = '  -0.21%  '

# Row 42
$ws.Range("D42").Value = "'95.15"
$ws.Range("E42").Value = '  -2.92%  '

# Row 43
$ws.Range("D43").Value = "'16.54"
$ws.Range("E43").Value = '  -0.31%  '

# Row 44
$ws.Range("D44").Value = "'0.0923"
$ws.Range("E44").Value = '  -2.93%  '

# Row 45
$ws.Range("E45").Value = '  -4.86%  '

# Row 46
$ws.Range("D46").Value = "'1.13"
$ws.Range("E46").Value = '  -4.48%  '

# Row 47
$ws.Range("D47").Value = "'7.27"
$ws.Range("E47").Value = '  +0.47%  '

# Row 48
$ws.Range("E48").Value = '  -2.14%  '

# Row 49
$ws.Range("E49").Value = '  -0.67%  '

# Row 50
$ws.Range("D50").Value = '2.205.51'
$ws.Range("E50").Value = '  -1.54%  '

# Row 51
$ws.Range("D51").Value = "'44.30"
